$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns D, K, L, M, N, O, P, Q, R, S, T for rows 2-28.
# (Columns A, B, C, E, F, G, H, I, J are constant and unchanged by this edit;
# the underlying edit is a re-shuffling/update of the weekly price rows.)
$rows = @(
    @(2, 44552, 'Castle Brite', 'Especial', 360, 20000, 21000, 20500, '$/caja 18 kilos', 'Región Metropolitana', 1139, 18),
    @(3, 44552, 'Castle Brite', 'Primera', 280, 18000, 19000, 18500, '$/caja 18 kilos', 'Región Metropolitana', 1028, 18),
    @(4, 44902, 'Castle Brite', 'Especial', 200, 25000, 26000, 25500, '$/caja 18 kilos', 'Región de O''Higgins', 1417, 18),
    @(5, 44902, 'Castle Brite', 'Primera', 240, 22000, 23000, 22500, '$/caja 18 kilos', 'Región de O''Higgins', 1250, 18),
    @(6, 44175, 'Castle Brite', 'Primera', 300, 21000, 22000, 21500, '$/caja 18 kilos', 'Región Metropolitana', 1194, 18),
    @(7, 44553, 'Modesto', 'Especial', 360, 23000, 24000, 23500, '$/caja 16 kilos', 'Región Metropolitana', 1469, 16),
    @(8, 44553, 'Modesto', 'Primera', 300, 21000, 22000, 21500, '$/caja 16 kilos', 'Región Metropolitana', 1344, 16),
    @(9, 44553, 'Modesto', 'Segunda', 240, 17000, 18000, 17500, '$/caja 16 kilos', 'Región Metropolitana', 1094, 16),
    @(10, 44189, 'Dina', 'Especial', 120, 23500, 24000, 23750, '$/caja 18 kilos', 'Región de O''Higgins', 1319, 18),
    @(11, 44189, 'Dina', 'Primera', 200, 21500, 22000, 21750, '$/caja 18 kilos', 'Región de O''Higgins', 1208, 18),
    @(12, 44546, 'Castle Brite', 'Especial', 300, 22500, 23000, 22750, '$/caja 18 kilos', 'Región Metropolitana', 1264, 18),
    @(13, 44546, 'Castle Brite', 'Primera', 300, 20500, 21000, 20750, '$/caja 18 kilos', 'Región Metropolitana', 1153, 18),
    @(14, 44160, 'Castle Brite', 'Primera', 240, 20500, 21000, 20750, '$/caja 15 kilos', 'Región Metropolitana', 1383, 15),
    @(15, 44161, 'Dina', 'Primera', 300, 20000, 20500, 20250, '$/caja 15 kilos', 'Región Metropolitana', 1350, 15),
    @(16, 44161, 'Dina', 'Segunda', 100, 18000, 18500, 18250, '$/caja 15 kilos', 'Región Metropolitana', 1217, 15),
    @(17, 44573, 'Modesto', 'Especial', 300, 20500, 21000, 20750, '$/caja 18 kilos', 'Región Metropolitana', 1153, 18),
    @(18, 44573, 'Modesto', 'Primera', 400, 17500, 18000, 17750, '$/caja 18 kilos', 'Región Metropolitana', 986, 18),
    @(19, 44566, 'Modesto', 'Especial', 100, 23000, 24000, 23500, '$/caja 18 kilos', 'Región de O''Higgins', 1306, 18),
    @(20, 44566, 'Modesto', 'Primera', 160, 21000, 22000, 21500, '$/caja 18 kilos', 'Región de O''Higgins', 1194, 18),
    @(21, 44895, 'Castle Brite', 'Primera', 200, 21000, 22000, 21500, '$/caja 16 kilos', 'Región Metropolitana', 1344, 16),
    @(22, 44559, 'Modesto', 'Especial', 400, 25000, 26000, 25500, '$/caja 18 kilos', 'Región de O''Higgins', 1417, 18),
    @(23, 44559, 'Modesto', 'Primera', 320, 22000, 23000, 22500, '$/caja 18 kilos', 'Región de O''Higgins', 1250, 18),
    @(24, 44545, 'Castle Brite', 'Especial', 340, 22500, 23000, 22750, '$/caja 18 kilos', 'Región de O''Higgins', 1264, 18),
    @(25, 44545, 'Castle Brite', 'Primera', 400, 20500, 21000, 20750, '$/caja 18 kilos', 'Región de O''Higgins', 1153, 18),
    @(26, 44545, 'Castle Brite', 'Segunda', 300, 15500, 16000, 15750, '$/caja 18 kilos', 'Región de O''Higgins', 875, 18),
    @(27, 44580, 'Modesto', 'Especial', 300, 22500, 23000, 22750, '$/caja 18 kilos', 'Región Metropolitana', 1264, 18),
    @(28, 44580, 'Modesto', 'Primera', 400, 19500, 20000, 19750, '$/caja 18 kilos', 'Región Metropolitana', 1097, 18)
)

foreach ($row in $rows) {
    $r  = $row[0]
    $ws.Cells.Item($r, 4).Value  = $row[1]   # D Fecha
    $ws.Cells.Item($r, 11).Value = $row[2]   # K Variedad
    $ws.Cells.Item($r, 12).Value = $row[3]   # L Calidad
    $ws.Cells.Item($r, 13).Value = $row[4]   # M Volumen
    $ws.Cells.Item($r, 14).Value = $row[5]   # N Precio minimo
    $ws.Cells.Item($r, 15).Value = $row[6]   # O Precio maximo
    $ws.Cells.Item($r, 16).Value = $row[7]   # P Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $row[8]   # Q Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $row[9]   # R Origen
    $ws.Cells.Item($r, 19).Value = $row[10]  # S Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $row[11]  # T Kg / unidad
}
